# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the per-language handback status sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-23 05:16:21"
$zhcn.Range("H2").Value = "2016-03-23 05:16:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-23 05:16:26"
$dede.Range("H2").Value = "2016-03-23 05:16:51"
